# ver 2.0.0 --> 2.0.1-beta1
# - Rename "Update Notes" --> "Release Notes" (title text boxes on slides 9 & 10)
# - Refresh cached date/slide-number placeholder text on the slide master and
#   all slide layouts (date field -> 16/05/2014, slide-number field -> <#>)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide master + every slide layout: refresh the cached placeholder text
#    for the "date" and "slide number" fields.
# ---------------------------------------------------------------------------
$lsq = [char]0x2039   # '<'  (single left-pointing angle quotation mark)
$rsq = [char]0x203a   # '>'  (single right-pointing angle quotation mark)
$newSlideNumText = "$lsq#$rsq"

$design = $ppt.ActivePresentation.Designs.Item(1)
$master = $design.SlideMaster

function Update-FieldPlaceholders($container) {
  for ($j = 1; $j -le $container.Shapes.Count; $j++) {
    $shp = $container.Shapes.Item($j)
    if ($shp.HasTextFrame) {
      if ($shp.Name -like "Datumsplatzhalter*") {
        $shp.TextFrame.TextRange.Text = "16/05/2014"
      } elseif ($shp.Name -like "Foliennummernplatzhalter*") {
        $shp.TextFrame.TextRange.Text = $newSlideNumText
      }
    }
  }
}

Update-FieldPlaceholders $master

$layouts = $master.CustomLayouts
for ($c = 1; $c -le $layouts.Count; $c++) {
  Update-FieldPlaceholders $layouts.Item($c)
}

# ---------------------------------------------------------------------------
# 2) Slides 9 & 10: rename "Update Notes.pptx" --> "Release Notes.pptx" in the
#    "Original in ..." caption text box, then let the autosize textbox regrow.
# ---------------------------------------------------------------------------
function Find-UpdateNotesShape($slideIndex) {
  $s = $p.Slides.Item($slideIndex)
  for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $shp = $s.Shapes.Item($j)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
      $full = $shp.TextFrame.TextRange.Text
      if ($full -like "*Update Notes.pptx*") {
        return $shp
      }
    }
  }
  return $null
}

# Slide 9: select from "Update" through the end and retype it in one go -
# the run holding "Update Notes.pptx" gets replaced by a single new run.
$shp9 = Find-UpdateNotesShape 9
$tr9 = $shp9.TextFrame.TextRange
$full9 = $tr9.Text
$idx9 = $full9.IndexOf("Update")
$tail9 = $tr9.Characters($idx9 + 1, $full9.Length - $idx9)
$tail9.Text = "Release Notes.pptx"
$shp9.Width = 267.1133070866142

# Slide 10: same word swap, plus the "in " run got re-touched separately,
# leaving it split off from "Original ".
$shp10 = Find-UpdateNotesShape 10
$tr10 = $shp10.TextFrame.TextRange
$full10 = $tr10.Text
$idx10 = $full10.IndexOf("Update")
$tail10 = $tr10.Characters($idx10 + 1, $full10.Length - $idx10)
$tail10.Text = "Release Notes.pptx"

$full10b = $tr10.Text
$idxIn = $full10b.IndexOf("in ")
$inWord = $tr10.Characters($idxIn + 1, 3)
$inWord.Text = "in "
$shp10.Width = 267.1133070866142

Write-Host "edit complete"
